$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Common")

$ws.Rows.Item(641).Insert()

$ws.Cells.Item(641, 1).Value = "cs"
$ws.Cells.Item(641, 2).Value = "MixtureInventory.list.total"
$ws.Cells.Item(641, 3).Value = "Počet mixů [{{data.total}}] ({{data.from}}-{{data.to}})"

$ws.Activate()
$ws.Range("B638").Select()
$excel.ActiveWindow.ScrollRow = 618
